$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns before existing "Tuin" column (G) -> new G,H blank; Tuin moves G->I, Score moves H->J
$ws.Columns("G:H").Insert()
# Insert two more blank columns before current Score column (J) -> J,K blank; Score moves J->L
$ws.Columns("J:K").Insert()

# Set the new header labels
$ws.Cells.Item(1, 7).Value = "Badkamers"
$ws.Cells.Item(1, 8).Value = "Bouwjaar"
$ws.Cells.Item(1, 10).Value = "Onderhoud Binnen"
$ws.Cells.Item(1, 11).Value = "Onderhoud Buiten"
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = " Gerard Doustraat 168 3, 1073 VZ Amsterdam Verkocht Width"
$ws.Cells.Item(2, 3).Value = 690000
$ws.Cells.Item(2, 4).Value = 84
$ws.Cells.Item(2, 5).Value = "C"
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 8).Value = 1876
$ws.Cells.Item(2, 9).Value = $false
$ws.Cells.Item(2, 10).Value = "Goed"
$ws.Cells.Item(2, 11).Value = "Goed"
$ws.Cells.Item(2, 12).Value = 0.7345454545454545
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = " Gerard Doustraat 2 3, 1072 VP Amsterdam Verkocht Width"
$ws.Cells.Item(3, 3).Value = 800000
$ws.Cells.Item(3, 4).Value = 87
$ws.Cells.Item(3, 5).Value = "B"
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = 1886
$ws.Cells.Item(3, 9).Value = $false
$ws.Cells.Item(3, 10).Value = "Goed"
$ws.Cells.Item(3, 11).Value = "Goed"
$ws.Cells.Item(3, 12).Value = 0.7272727272727272
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = " Quellijnstraat 56 3, 1072 XT Amsterdam Verkocht Width"
$ws.Cells.Item(4, 3).Value = 649000
$ws.Cells.Item(4, 4).Value = 67
$ws.Cells.Item(4, 5).Value = "A"
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 2
$ws.Cells.Item(4, 8).Value = 1906
$ws.Cells.Item(4, 9).Value = $false
$ws.Cells.Item(4, 10).Value = "Uitstekend"
$ws.Cells.Item(4, 11).Value = "Goed"
$ws.Cells.Item(4, 12).Value = 0.6920833333333334
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = " Gerard Doustraat 3 A3, 1072 VH Amsterdam Verkocht Width"
$ws.Cells.Item(5, 3).Value = 640000
$ws.Cells.Item(5, 4).Value = 79
$ws.Cells.Item(5, 5).Value = "C"
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 2
$ws.Cells.Item(5, 8).Value = 1880
$ws.Cells.Item(5, 9).Value = $false
$ws.Cells.Item(5, 10).Value = "Goed"
$ws.Cells.Item(5, 11).Value = "Goed"
$ws.Cells.Item(5, 12).Value = 0.6912121212121213
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = " Gerard Doustraat 234 2, 1073 XC Amsterdam Verkocht Width"
$ws.Cells.Item(6, 3).Value = 475000
$ws.Cells.Item(6, 4).Value = 58
$ws.Cells.Item(6, 5).Value = "C"
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 1
$ws.Cells.Item(6, 8).Value = 1892
$ws.Cells.Item(6, 9).Value = $false
$ws.Cells.Item(6, 10).Value = "Goed"
$ws.Cells.Item(6, 11).Value = "Goed"
$ws.Cells.Item(6, 12).Value = 0.6878787878787879
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = " Eerste Sweelinckstraat 5 3, 1073 CK Amsterdam Verkocht Width"
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 70
$ws.Cells.Item(7, 5).Value = "A"
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 1886
$ws.Cells.Item(7, 9).Value = $false
$ws.Cells.Item(7, 10).Value = "Goed"
$ws.Cells.Item(7, 11).Value = "Goed"
$ws.Cells.Item(7, 12).Value = 0.6833333333333333
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = " Quellijnstraat 19 A, 1072 XM Amsterdam Verkocht Width"
$ws.Cells.Item(8, 3).Value = 725000
$ws.Cells.Item(8, 4).Value = 77
$ws.Cells.Item(8, 5).Value = "C"
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 1
$ws.Cells.Item(8, 8).Value = 1879
$ws.Cells.Item(8, 9).Value = $false
$ws.Cells.Item(8, 10).Value = "Goed"
$ws.Cells.Item(8, 11).Value = "Goed"
$ws.Cells.Item(8, 12).Value = 0.6599621212121213
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = " Gerard Doustraat 192 2, 1073 XA Amsterdam Verkocht Width"
$ws.Cells.Item(9, 3).Value = 450000
$ws.Cells.Item(9, 4).Value = 50
$ws.Cells.Item(9, 5).Value = "D"
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 1
$ws.Cells.Item(9, 8).Value = 1906
$ws.Cells.Item(9, 9).Value = $false
$ws.Cells.Item(9, 10).Value = "Goed"
$ws.Cells.Item(9, 11).Value = "Goed"
$ws.Cells.Item(9, 12).Value = 0.6584848484848485
$ws.Cells.Item(10, 1).Value = 9
$ws.Cells.Item(10, 2).Value = " Quellijnstraat 17 3, 1072 XM Amsterdam Verkocht Width"
$ws.Cells.Item(10, 3).Value = 675000
$ws.Cells.Item(10, 4).Value = 76
$ws.Cells.Item(10, 5).Value = "C"
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 1
$ws.Cells.Item(10, 8).Value = 1879
$ws.Cells.Item(10, 9).Value = $false
$ws.Cells.Item(10, 10).Value = "Goed"
$ws.Cells.Item(10, 11).Value = "Goed"
$ws.Cells.Item(10, 12).Value = 0.6566287878787879
$ws.Cells.Item(11, 1).Value = 10
$ws.Cells.Item(11, 2).Value = " Quellijnstraat 41 B, 1072 XP Amsterdam Verkocht Width"
$ws.Cells.Item(11, 3).Value = 785000
$ws.Cells.Item(11, 4).Value = 74
$ws.Cells.Item(11, 5).Value = "A"
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 1
$ws.Cells.Item(11, 8).Value = 1906
$ws.Cells.Item(11, 9).Value = $false
$ws.Cells.Item(11, 10).Value = "Goed"
$ws.Cells.Item(11, 11).Value = "Goed"
$ws.Cells.Item(11, 12).Value = 0.6554166666666668
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = " Gerard Doustraat 147 E, 1073 VV Amsterdam Verkocht Width"
$ws.Cells.Item(12, 3).Value = 800000
$ws.Cells.Item(12, 4).Value = 87
$ws.Cells.Item(12, 5).Value = "Unknown"
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 1
$ws.Cells.Item(12, 8).Value = 2010
$ws.Cells.Item(12, 9).Value = $false
$ws.Cells.Item(12, 10).Value = "Goed"
$ws.Cells.Item(12, 11).Value = "Goed"
$ws.Cells.Item(12, 12).Value = 0.6549999999999999
$ws.Cells.Item(13, 1).Value = 12
$ws.Cells.Item(13, 2).Value = " Quellijnstraat 37 B, 1072 XP Amsterdam Verkocht Width"
$ws.Cells.Item(13, 3).Value = 725000
$ws.Cells.Item(13, 4).Value = 74
$ws.Cells.Item(13, 5).Value = "B"
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 1
$ws.Cells.Item(13, 8).Value = 1879
$ws.Cells.Item(13, 9).Value = $false
$ws.Cells.Item(13, 10).Value = "Goed"
$ws.Cells.Item(13, 11).Value = "Goed"
$ws.Cells.Item(13, 12).Value = 0.652689393939394
$ws.Cells.Item(14, 1).Value = 13
$ws.Cells.Item(14, 2).Value = " Gerard Doustraat 82 3, 1072 VW Amsterdam Verkocht Width"
$ws.Cells.Item(14, 3).Value = 565000
$ws.Cells.Item(14, 4).Value = 79
$ws.Cells.Item(14, 5).Value = "G"
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 1876
$ws.Cells.Item(14, 9).Value = $false
$ws.Cells.Item(14, 10).Value = "Goed"
$ws.Cells.Item(14, 11).Value = "Goed"
$ws.Cells.Item(14, 12).Value = 0.646969696969697
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = " Eerste Jan van der Heijdenstraat 32 3, 1072 TV Amsterdam Verkocht Width"
$ws.Cells.Item(15, 3).Value = 750000
$ws.Cells.Item(15, 4).Value = 78
$ws.Cells.Item(15, 5).Value = "A"
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(15, 8).Value = 1906
$ws.Cells.Item(15, 9).Value = $false
$ws.Cells.Item(15, 10).Value = "Goed"
$ws.Cells.Item(15, 11).Value = "Goed"
$ws.Cells.Item(15, 12).Value = 0.640625
$ws.Cells.Item(16, 1).Value = 15
$ws.Cells.Item(16, 2).Value = " Eerste Jan van der Heijdenstraat 42 1, 1072 TV Amsterdam Verkocht Width"
$ws.Cells.Item(16, 3).Value = 675000
$ws.Cells.Item(16, 4).Value = 77
$ws.Cells.Item(16, 5).Value = "A"
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 1
$ws.Cells.Item(16, 8).Value = 1906
$ws.Cells.Item(16, 9).Value = $false
$ws.Cells.Item(16, 10).Value = "Goed"
$ws.Cells.Item(16, 11).Value = "Goed"
$ws.Cells.Item(16, 12).Value = 0.6372916666666668
